$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the filename value in B8: "Laos24" -> "Laos23"
$ws.Range("B8").Value = "Laos23"

# 2. Give B1 (empty header cell) the same formatting as B4:B7 (right aligned,
#    bordered, bold font) by copying the format from B4.
$ws.Range("B4").Copy() | Out-Null
$ws.Range("B1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# 3. Give B8 (the filename cell) the same formatting as B2/B3 (left aligned,
#    bordered font) by copying the format from B2.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# 4. Row heights for all data rows grew slightly (18.75 -> 19.5)
for ($r = 1; $r -le 8; $r++) {
    $ws.Rows.Item($r).RowHeight = 19.5
}
